$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for I ("I0") and J ("IF") columns, rows 2-28
$data = @(
    @(1, 3),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 5),
    @(3, 4),
    @(9, 9),
    @(1, 2),
    @(5, 5),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
